# Weekly update: insert a new price observation for "Apio" at
# Terminal Hortofrutícola Agro Chillán, pushing the existing rows
# (122:141) down by one and appending what was the last row as row 142.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 122:141 down to 123:142 (Excel re-flows the existing values
# automatically, exactly like Insert does in the real UI).
$ws.Rows.Item(122).Insert()

# Populate the newly-opened row 122 with this week's data point.
$ws.Cells.Item(122, 1).Value  = 7
$ws.Cells.Item(122, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(122, 3).Value  = "Ñuble"
$ws.Cells.Item(122, 4).Value  = 44505
$ws.Cells.Item(122, 5).Value  = 16
$ws.Cells.Item(122, 6).Value  = 100112017
$ws.Cells.Item(122, 7).Value  = "Apio"
$ws.Cells.Item(122, 8).Value  = "Americana (o)"
$ws.Cells.Item(122, 9).Value  = "Primera"
$ws.Cells.Item(122, 10).Value = 100
$ws.Cells.Item(122, 11).Value = 8000
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = 8500
$ws.Cells.Item(122, 14).Value = "$/docena de matas"
$ws.Cells.Item(122, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(122, 16).Value = 1417
$ws.Cells.Item(122, 17).Value = 6
$ws.Cells.Item(122, 18).Value = "Hortaliza"
